$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add F1 header, matching style (bold, centered, bordered) of other header cells
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data rows (row, A-serial, B, C, D, E, F)
$rows = @(
    @(2,  45685.64698888889, 1014.8, 12.86, 3.465544564383369, "10-15", "Duża Gra"),
    @(3,  45685.64958379629, 1239,   12.76, 3.530149357659476, "10-15", "Duża Gra"),
    @(4,  45685.66118449074, 2241.3, 12.77, 3.555028438568115, "10-15", "Duża Gra"),
    @(5,  45685.64698541666, 1014.5, 8.970000000000001, 3.199397053037372, "5-10", "Duża Gra"),
    @(6,  45685.64958032408, 1238.7, 8.92,  3.251099143709456, "5-10", "Duża Gra"),
    @(7,  45685.66118217593, 2241.1, 9.9,   3.263054694448197, "5-10", "Duża Gra"),
    @(8,  45685.67135,       3119.6, 13.47, 3.462289776120867, "10-15", "Mała Gra"),
    @(9,  45685.67190324074, 3167.4, 12.6,  3.167921134403773, "10-15", "Mała Gra"),
    @(10, 45685.68490787037, 4291,   11.67, 2.918576104300363, "10-15", "Mała Gra"),
    @(11, 45685.67134652778, 3119.3, 9.15,  2.776583075523376, "5-10", "Mała Gra"),
    @(12, 45685.67189976852, 3167.1, 8.779999999999999, 2.813431944165912, "5-10", "Mała Gra"),
    @(13, 45685.68056759259, 3916,   9.720000000000001, 2.817018202372959, "5-10", "Mała Gra")
)

$firstDateCell = $true
foreach ($r in $rows) {
    $rowNum = $r[0]
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r[1]
    if ($firstDateCell) {
        # Apply lowercase format first, then uppercase. This registers both
        # numFmt entries (164 lowercase, 165 uppercase) while only the
        # uppercase one ends up referenced by a cellXf.
        $cellA.NumberFormat = "yyyy-mm-dd h:mm:ss"
        $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $firstDateCell = $false
    } else {
        $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }

    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}
